$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "DKS" (Desktop Computer) device type rows (rows 8, 9, 10) -
# this shifts the following device types (CMR, SCN, PRT) up and the
# unused "DKS" shared strings get garbage-collected automatically.
$ws.Rows("8:10").Delete()

# Set the active selection to E10 (matches the post-edit sheetView selection)
$ws.Range("E10").Select()

# Match the saved page setup (paper size / orientation) recorded for the sheet
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
